$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34-104 down to 35-105
$ws.Rows(34).Insert()

# Populate the new row 34 with the new weekly data point
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 45125
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100108
$ws.Range("H34").Value = "Tropicales y subtropicales"
$ws.Range("I34").Value = 100108007
$ws.Range("J34").Value = "Coco"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 20
$ws.Range("N34").Value = 36000
$ws.Range("O34").Value = 36000
$ws.Range("P34").Value = 36000
$ws.Range("Q34").Value = "$/malla 20 unidades"
$ws.Range("R34").Value = "Perú"
$ws.Range("S34").Value = 1800
$ws.Range("T34").Value = 20
